$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("I8").Value = "fc33f22e-ddd0-4d1b-be08-49b763b92fac.md"
$ws2.Range("I8").Style = "Hyperlink"
